# Update odds values in rows 2, 5, 6, 7 as per the source update (Atualizando o arquivo XLSX)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
${ws}.Range("G2").Value = 2.27
${ws}.Range("H2").Value = 2.67
${ws}.Range("I2").Value = 3.75
${ws}.Range("J2").Value = 3.1
${ws}.Range("L2").Value = 4.65
${ws}.Range("N2").Value = 4.25
${ws}.Range("R2").Value = 1.29
${ws}.Range("U2").Value = 1.75
${ws}.Range("V2").Value = 1.98
${ws}.Range("W2").Value = 2.55
${ws}.Range("X2").Value = 1.45
${ws}.Range("Y2").Value = 4.75
${ws}.Range("Z2").Value = 8.75
${ws}.Range("AA2").Value = 11
${ws}.Range("AB2").Value = 23
${ws}.Range("AC2").Value = 29
${ws}.Range("AE2").Value = 4.25
${ws}.Range("AF2").Value = 5.9
${ws}.Range("AI2").Value = 6.6
${ws}.Range("AJ2").Value = 17
${ws}.Range("AK2").Value = 15.5
${ws}.Range("AL2").Value = 65
${ws}.Range("AM2").Value = 55

# Row 5
${ws}.Range("G5").Value = 1.33
${ws}.Range("J5").Value = 1.73
${ws}.Range("K5").Value = 2.63
${ws}.Range("S5").Value = 2
${ws}.Range("T5").Value = 1.73

# Row 6
${ws}.Range("G6").Value = 1.27
${ws}.Range("J6").Value = 1.67
${ws}.Range("K6").Value = 2.63
${ws}.Range("N6").Value = 12
${ws}.Range("O6").Value = 1.13
${ws}.Range("P6").Value = 5.5

# Row 7
${ws}.Range("G7").Value = 2.12
${ws}.Range("H7").Value = 3.25
${ws}.Range("I7").Value = 3.3
${ws}.Range("J7").Value = 2.67
${ws}.Range("K7").Value = 2.12
${ws}.Range("L7").Value = 3.8
${ws}.Range("M7").Value = 1.08
${ws}.Range("N7").Value = 6.6
${ws}.Range("O7").Value = 1.37
${ws}.Range("P7").Value = 2.85
${ws}.Range("Q7").Value = 2.1
${ws}.Range("R7").Value = 1.65
${ws}.Range("S7").Value = 3.6
${ws}.Range("T7").Value = 1.25
${ws}.Range("U7").Value = 1.4
${ws}.Range("V7").Value = 2.72
${ws}.Range("W7").Value = 1.88
${ws}.Range("X7").Value = 1.82
${ws}.Range("Y7").Value = 6.7
${ws}.Range("Z7").Value = 9.5
${ws}.Range("AA7").Value = 9
${ws}.Range("AB7").Value = 19.5
${ws}.Range("AC7").Value = 18
${ws}.Range("AD7").Value = 32
${ws}.Range("AE7").Value = 6.6
${ws}.Range("AF7").Value = 6.3
${ws}.Range("AG7").Value = 15.5
${ws}.Range("AI7").Value = 8.75
${ws}.Range("AJ7").Value = 16.5
${ws}.Range("AK7").Value = 11.75
${ws}.Range("AL7").Value = 45
${ws}.Range("AM7").Value = 32

